$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts the assigned string
# into a numeric value (losing formatting / introducing float noise).
$textForceCells = @("D6","D7","D10","D13","D18","D20","D21","D23","D26","D30","D34","D35","D37","D38","D42","D43","D44","D45","D47","D48","D49","D50")
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "97.616.85"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "3.725.15"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +13.51%  "
$ws.Range("D6").Value = "238.45"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "658.06"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  +5.31%  "
$ws.Range("E9").Value = "  +4.43%  "
$ws.Range("D10").Value = "1.00"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").Value = "3.727.07"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("E12").Value = "  +17.64%  "
$ws.Range("D13").Value = "44.88"
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "4.425.93"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "97.419.32"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "9.25"
$ws.Range("E18").Value = "  +2.88%  "
$ws.Range("D19").Value = "3.728.27"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").Value = "13.13"
$ws.Range("E20").Value = "  +2.75%  "
$ws.Range("D21").Value = "18.87"
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("E22").Value = "  +2.65%  "
$ws.Range("D23").Value = "529.78"
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  +10.41%  "
$ws.Range("D26").Value = "117.92"
$ws.Range("E26").Value = "  +15.59%  "
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("E28").Value = "  +25.93%  "
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").Value = "12.82"
$ws.Range("E30").Value = "  +2.50%  "
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  +2.99%  "
$ws.Range("D34").Value = "1.84"
$ws.Range("E34").Value = "  -2.90%  "
$ws.Range("D35").Value = "33.18"
$ws.Range("E35").Value = "  +1.51%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").Value = "0.599"
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("D38").Value = "642.25"
$ws.Range("E38").Value = "  -3.12%  "
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  +5.50%  "
$ws.Range("D42").Value = "6.84"
$ws.Range("E42").Value = "  -3.68%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "41.21"
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "0.495"
$ws.Range("E44").Value = "  +9.66%  "
$ws.Range("D45").Value = "2.01"
$ws.Range("E45").Value = "  +2.01%  "
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").Value = "0.0459"
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("D48").Value = "2.40"
$ws.Range("E48").Value = "  +2.97%  "
$ws.Range("D49").Value = "8.80"
$ws.Range("E49").Value = "  +2.66%  "
$ws.Range("D50").Value = "23.66"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("E51").Value = "  +4.81%  "
